$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1043256666666667
$ws.Range("H2").Value = 0.312977
$ws.Range("I2").Value = 0.02547563162231953
$ws.Range("J2").Value = 0.02547563162231953
$ws.Range("M2").Value = 1.484826
$ws.Range("N2").Value = 4.454478
$ws.Range("O2").Value = 0.06049021884829667
$ws.Range("P2").Value = 0.06049021884829667
$ws.Range("Q2").Value = 0.154905462334
$ws.Range("R2").Value = 1.394149161006
$ws.Range("S2").Value = 0.001541026532132696
$ws.Range("T2").Value = 0.001541026532132696

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1043256666666667
$ws.Range("H3").Value = 0.312977
$ws.Range("I3").Value = 0.02547563162231953
$ws.Range("J3").Value = 0.02547563162231953
$ws.Range("O3").Value = 0.5859425360316464
$ws.Range("P3").Value = 0.5859425360316464
$ws.Range("Q3").Value = 1.500502084027333
$ws.Range("R3").Value = 13.504518756246
$ws.Range("S3").Value = 0.01492725619978991
$ws.Range("T3").Value = 0.01492725619978991

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1043256666666667
$ws.Range("H4").Value = 0.312977
$ws.Range("I4").Value = 0.02547563162231953
$ws.Range("J4").Value = 0.02547563162231953
$ws.Range("M4").Value = 8.653369666666666
$ws.Range("N4").Value = 25.960109
$ws.Range("O4").Value = 0.3525289999716321
$ws.Range("P4").Value = 0.3525289999716321
$ws.Range("Q4").Value = 0.9027685593881111
$ws.Range("R4").Value = 8.124917034493
$ws.Range("S4").Value = 0.008980898939461992
$ws.Range("T4").Value = 0.008980898939461992

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1043256666666667
$ws.Range("H5").Value = 0.312977
$ws.Range("I5").Value = 0.02547563162231953
$ws.Range("J5").Value = 0.02547563162231953
$ws.Range("M5").Value = 0.02548533333333333
$ws.Range("N5").Value = 0.076456
$ws.Range("O5").Value = 0.001038245148424882
$ws.Range("P5").Value = 0.001038245148424882
$ws.Range("Q5").Value = 0.002658774390222222
$ws.Range("R5").Value = 0.023928969512
$ws.Range("S5").Value = 0.00002644995093493275
$ws.Range("T5").Value = 0.00002644995093493275

# Row 6
$ws.Range("I6").Value = 0.9745243683776804
$ws.Range("J6").Value = 0.9745243683776804
$ws.Range("M6").Value = 1.484826
$ws.Range("N6").Value = 4.454478
$ws.Range("O6").Value = 0.06049021884829667
$ws.Range("P6").Value = 0.06049021884829667
$ws.Range("Q6").Value = 5.925629247481999
$ws.Range("R6").Value = 53.33066322733799
$ws.Range("S6").Value = 0.05894919231616397
$ws.Range("T6").Value = 0.05894919231616397

# Row 7
$ws.Range("I7").Value = 0.9745243683776804
$ws.Range("J7").Value = 0.9745243683776804
$ws.Range("O7").Value = 0.5859425360316464
$ws.Range("P7").Value = 0.5859425360316464
$ws.Range("Q7").Value = 57.39900259842866
$ws.Range("R7").Value = 516.591023385858
$ws.Range("S7").Value = 0.5710152798318564
$ws.Range("T7").Value = 0.5710152798318564

# Row 8
$ws.Range("I8").Value = 0.9745243683776804
$ws.Range("J8").Value = 0.9745243683776804
$ws.Range("M8").Value = 8.653369666666666
$ws.Range("N8").Value = 25.960109
$ws.Range("O8").Value = 0.3525289999716321
$ws.Range("P8").Value = 0.3525289999716321
$ws.Range("Q8").Value = 34.53378401649322
$ws.Range("R8").Value = 310.804056148439
$ws.Range("S8").Value = 0.3435481010321701
$ws.Range("T8").Value = 0.3435481010321701

# Row 9
$ws.Range("I9").Value = 0.9745243683776804
$ws.Range("J9").Value = 0.9745243683776804
$ws.Range("M9").Value = 0.02548533333333333
$ws.Range("N9").Value = 0.076456
$ws.Range("O9").Value = 0.001038245148424882
$ws.Range("P9").Value = 0.001038245148424882
$ws.Range("Q9").Value = 0.1017066219084444
$ws.Range("R9").Value = 0.9153595971759999
$ws.Range("S9").Value = 0.001011795197489949
$ws.Range("T9").Value = 0.001011795197489949
